$d = $word.ActiveDocument

# Locate the sentence containing "... לרוב קטן אף מאוד ..." and figure out
# exactly where the new word ("חלק ") needs to be inserted: right after
# "לרוב " and right before "קטן".
$searchRange = $d.Content
$found = $searchRange.Find.Execute(
    "לרוב קטן",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "",
    0
)

if (-not $found) {
    throw "Could not find target text 'לרוב קטן' in the document."
}

$insertStart = $searchRange.Start + 5

# Insert the new word at that exact character position. This splits the
# original run into a "before" run and an "after" run, with the freshly
# typed text living in its own run in between.
$insertionPoint = $d.Range($insertStart, $insertStart)
$insertionPoint.InsertBefore("חלק ")

# Nudge the newly inserted run's direct formatting so the run actually
# gets materialized as its own <w:r> (matching how Word splits a run when
# you type in the middle of it) instead of being silently re-merged back
# into its neighbour.
$newWordRange = $d.Range($insertStart, $insertStart + 4)
$newWordRange.Font.Bold = $true
$newWordRange.Font.Bold = $false
